# Add a "Save" column (H) to the s_vals sheet, mirroring the header style
# used by the existing columns (B1:G1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the neighboring header cell (G1) onto the new
# header cell (H1) so it picks up the same bold/border/alignment style.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("H1").Value = "Save"

$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 1
